$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("保險")

# Column B: company (header + data), written top-to-bottom so new shared
# strings are appended in the same order the original edit produced them.
$ws.Range("B1").Value = "company"
$ws.Range("B2").Value = "南山人壽"
$ws.Range("B3").Value = "南山人壽"
$ws.Range("B4").Value = "南山人壽"
$ws.Range("B5").Value = "全球人壽"
$ws.Range("B6").Value = "全球人壽"
$ws.Range("B7").Value = "全球人壽"
$ws.Range("B8").Value = "國寶人壽"
$ws.Range("B9").Value = "南山人壽"
$ws.Range("B10").Value = "國寶人壽"
$ws.Range("B11").Value = "國寶人壽"
$ws.Range("B12").Value = "南山人壽"
$ws.Range("B13").Value = "國寶人壽"
$ws.Range("B14").Value = "國寶人壽"
$ws.Range("B15").Value = "國寶人壽"
$ws.Range("B16").Value = "中國人壽"
$ws.Range("B17").Value = "中國人壽"
$ws.Range("B18").Value = "中國人壽"
$ws.Range("B19").Value = "中國人壽"
$ws.Range("B20").Value = "幸福人壽"
$ws.Range("B21").Value = "中國人壽"
$ws.Range("B22").Value = "幸福人壽"
$ws.Range("B23").Value = "安聯人壽"
$ws.Range("B24").Value = "幸福人壽"
$ws.Range("B25").Value = "安聯人壽"
$ws.Range("B26").Value = "幸福人壽"
$ws.Range("B27").Value = "安聯人壽"
$ws.Range("B28").Value = "幸福人壽"
$ws.Range("B29").Value = "富邦人壽"
$ws.Range("B30").Value = "富邦人壽"

# Column C: name (header + data)
$ws.Range("C1").Value = "name"
$ws.Range("C2").Value = "新20年限期繳費增值分紅终身保險"
$ws.Range("C3").Value = "增值分紅終身壽險"
$ws.Range("C4").Value = "增值分紅終身"
$ws.Range("C5").Value = "inj增額終身B型"
$ws.Range("C6").Value = "還本終身"
$ws.Range("C7").Value = "增額終身B型"
$ws.Range("C8").Value = "長青增額终身險"
$ws.Range("C9").Value = "家庭防癌險"
$ws.Range("C10").Value = "永泰終身"
$ws.Range("C11").Value = "永泰終身"
$ws.Range("C12").Value = "家庭防癌險"
$ws.Range("C13").Value = "永泰终身"
$ws.Range("C14").Value = "永泰終身險"
$ws.Range("C15").Value = "永泰終身險"
$ws.Range("C16").Value = "金享受"
$ws.Range("C17").Value = "金享受"
$ws.Range("C18").Value = "金享受"
$ws.Range("C19").Value = "鑫幸福終身壽險"
$ws.Range("C20").Value = "新防癌終身ICG"
$ws.Range("C21").Value = "鑫幸福終身壽險"
$ws.Range("C22").Value = "新防癌終身ICG"
$ws.Range("C23").Value = "五年定期重大疾病"
$ws.Range("C24").Value = "新防癌終身ICG"
$ws.Range("C25").Value = "五年定期重大疾病"
$ws.Range("C26").Value = "新防癌终身IOG"
$ws.Range("C27").Value = "五年定期重大疾病"
$ws.Range("C28").Value = "新防癌終身IOG"
$ws.Range("C29").Value = "幸福誠保意外險"
$ws.Range("C30").Value = "幸福誠保意外險"

# Column D: owner (header + data)
$ws.Range("D1").Value = "owner"
$ws.Range("D2").Value = "林德福"
$ws.Range("D3").Value = "林德福"
$ws.Range("D4").Value = "林德福"
$ws.Range("D5").Value = "林德福"
$ws.Range("D6").Value = "曾月桂"
$ws.Range("D7").Value = "曾月桂"
$ws.Range("D8").Value = "曾月桂"
$ws.Range("D9").Value = "林德福"
$ws.Range("D10").Value = "林德福"
$ws.Range("D11").Value = "曾月桂"
$ws.Range("D12").Value = "林德福"
$ws.Range("D13").Value = "曾月桂"
$ws.Range("D14").Value = "曾月麁"
$ws.Range("D15").Value = "曾月桂"
$ws.Range("D16").Value = "曾月桂"
$ws.Range("D17").Value = "曾月桂"
$ws.Range("D18").Value = "曾月桂"
$ws.Range("D19").Value = "曾月桂"
$ws.Range("D20").Value = "曾月桂"
$ws.Range("D21").Value = "林德福"
$ws.Range("D22").Value = "林德福"
$ws.Range("D23").Value = "曾月桂"
$ws.Range("D24").Value = "曾月桂"
$ws.Range("D25").Value = "曾月桂"
$ws.Range("D26").Value = "曾月桂"
$ws.Range("D27").Value = "曾月桂"
$ws.Range("D28").Value = "曾月桂"
$ws.Range("D29").Value = "曾月桂"
$ws.Range("D30").Value = "曾月桂"

# Column E: property_category (header + data)
$ws.Range("E1").Value = "property_category"
$ws.Range("E2").Value = "insurance"
$ws.Range("E3").Value = "insurance"
$ws.Range("E4").Value = "insurance"
$ws.Range("E5").Value = "insurance"
$ws.Range("E6").Value = "insurance"
$ws.Range("E7").Value = "insurance"
$ws.Range("E8").Value = "insurance"
$ws.Range("E9").Value = "insurance"
$ws.Range("E10").Value = "insurance"
$ws.Range("E11").Value = "insurance"
$ws.Range("E12").Value = "insurance"
$ws.Range("E13").Value = "insurance"
$ws.Range("E14").Value = "insurance"
$ws.Range("E15").Value = "insurance"
$ws.Range("E16").Value = "insurance"
$ws.Range("E17").Value = "insurance"
$ws.Range("E18").Value = "insurance"
$ws.Range("E19").Value = "insurance"
$ws.Range("E20").Value = "insurance"
$ws.Range("E21").Value = "insurance"
$ws.Range("E22").Value = "insurance"
$ws.Range("E23").Value = "insurance"
$ws.Range("E24").Value = "insurance"
$ws.Range("E25").Value = "insurance"
$ws.Range("E26").Value = "insurance"
$ws.Range("E27").Value = "insurance"
$ws.Range("E28").Value = "insurance"
$ws.Range("E29").Value = "insurance"
$ws.Range("E30").Value = "insurance"

# Column F: category (header + data)
$ws.Range("F1").Value = "category"
$ws.Range("F2").Value = "normal"
$ws.Range("F3").Value = "normal"
$ws.Range("F4").Value = "normal"
$ws.Range("F5").Value = "normal"
$ws.Range("F6").Value = "normal"
$ws.Range("F7").Value = "normal"
$ws.Range("F8").Value = "normal"
$ws.Range("F9").Value = "normal"
$ws.Range("F10").Value = "normal"
$ws.Range("F11").Value = "normal"
$ws.Range("F12").Value = "normal"
$ws.Range("F13").Value = "normal"
$ws.Range("F14").Value = "normal"
$ws.Range("F15").Value = "normal"
$ws.Range("F16").Value = "normal"
$ws.Range("F17").Value = "normal"
$ws.Range("F18").Value = "normal"
$ws.Range("F19").Value = "normal"
$ws.Range("F20").Value = "normal"
$ws.Range("F21").Value = "normal"
$ws.Range("F22").Value = "normal"
$ws.Range("F23").Value = "normal"
$ws.Range("F24").Value = "normal"
$ws.Range("F25").Value = "normal"
$ws.Range("F26").Value = "normal"
$ws.Range("F27").Value = "normal"
$ws.Range("F28").Value = "normal"
$ws.Range("F29").Value = "normal"
$ws.Range("F30").Value = "normal"

# Column G: date (header + data)
$ws.Range("G1").Value = "date"
$ws.Range("G2").Value = "2011-11-28"
$ws.Range("G3").Value = "2011-11-28"
$ws.Range("G4").Value = "2011-11-28"
$ws.Range("G5").Value = "2011-11-28"
$ws.Range("G6").Value = "2011-11-28"
$ws.Range("G7").Value = "2011-11-28"
$ws.Range("G8").Value = "2011-11-28"
$ws.Range("G9").Value = "2011-11-28"
$ws.Range("G10").Value = "2011-11-28"
$ws.Range("G11").Value = "2011-11-28"
$ws.Range("G12").Value = "2011-11-28"
$ws.Range("G13").Value = "2011-11-28"
$ws.Range("G14").Value = "2011-11-28"
$ws.Range("G15").Value = "2011-11-28"
$ws.Range("G16").Value = "2011-11-28"
$ws.Range("G17").Value = "2011-11-28"
$ws.Range("G18").Value = "2011-11-28"
$ws.Range("G19").Value = "2011-11-28"
$ws.Range("G20").Value = "2011-11-28"
$ws.Range("G21").Value = "2011-11-28"
$ws.Range("G22").Value = "2011-11-28"
$ws.Range("G23").Value = "2011-11-28"
$ws.Range("G24").Value = "2011-11-28"
$ws.Range("G25").Value = "2011-11-28"
$ws.Range("G26").Value = "2011-11-28"
$ws.Range("G27").Value = "2011-11-28"
$ws.Range("G28").Value = "2011-11-28"
$ws.Range("G29").Value = "2011-11-28"
$ws.Range("G30").Value = "2011-11-28"

# Column H: legislator_name (header + data)
$ws.Range("H1").Value = "legislator_name"
$ws.Range("H2").Value = "林德福"
$ws.Range("H3").Value = "林德福"
$ws.Range("H4").Value = "林德福"
$ws.Range("H5").Value = "林德福"
$ws.Range("H6").Value = "林德福"
$ws.Range("H7").Value = "林德福"
$ws.Range("H8").Value = "林德福"
$ws.Range("H9").Value = "林德福"
$ws.Range("H10").Value = "林德福"
$ws.Range("H11").Value = "林德福"
$ws.Range("H12").Value = "林德福"
$ws.Range("H13").Value = "林德福"
$ws.Range("H14").Value = "林德福"
$ws.Range("H15").Value = "林德福"
$ws.Range("H16").Value = "林德福"
$ws.Range("H17").Value = "林德福"
$ws.Range("H18").Value = "林德福"
$ws.Range("H19").Value = "林德福"
$ws.Range("H20").Value = "林德福"
$ws.Range("H21").Value = "林德福"
$ws.Range("H22").Value = "林德福"
$ws.Range("H23").Value = "林德福"
$ws.Range("H24").Value = "林德福"
$ws.Range("H25").Value = "林德福"
$ws.Range("H26").Value = "林德福"
$ws.Range("H27").Value = "林德福"
$ws.Range("H28").Value = "林德福"
$ws.Range("H29").Value = "林德福"
$ws.Range("H30").Value = "林德福"

# Column I: legislator_id (header + data)
$ws.Range("I1").Value = "legislator_id"
$ws.Range("I2").Value = 908
$ws.Range("I3").Value = 908
$ws.Range("I4").Value = 908
$ws.Range("I5").Value = 908
$ws.Range("I6").Value = 908
$ws.Range("I7").Value = 908
$ws.Range("I8").Value = 908
$ws.Range("I9").Value = 908
$ws.Range("I10").Value = 908
$ws.Range("I11").Value = 908
$ws.Range("I12").Value = 908
$ws.Range("I13").Value = 908
$ws.Range("I14").Value = 908
$ws.Range("I15").Value = 908
$ws.Range("I16").Value = 908
$ws.Range("I17").Value = 908
$ws.Range("I18").Value = 908
$ws.Range("I19").Value = 908
$ws.Range("I20").Value = 908
$ws.Range("I21").Value = 908
$ws.Range("I22").Value = 908
$ws.Range("I23").Value = 908
$ws.Range("I24").Value = 908
$ws.Range("I25").Value = 908
$ws.Range("I26").Value = 908
$ws.Range("I27").Value = 908
$ws.Range("I28").Value = 908
$ws.Range("I29").Value = 908
$ws.Range("I30").Value = 908

# Column J: source_file (header + data)
$ws.Range("J1").Value = "source_file"
$ws.Range("J2").Value = "tmp2e4a1"
$ws.Range("J3").Value = "tmp2e4a1"
$ws.Range("J4").Value = "tmp2e4a1"
$ws.Range("J5").Value = "tmp2e4a1"
$ws.Range("J6").Value = "tmp2e4a1"
$ws.Range("J7").Value = "tmp2e4a1"
$ws.Range("J8").Value = "tmp2e4a1"
$ws.Range("J9").Value = "tmp2e4a1"
$ws.Range("J10").Value = "tmp2e4a1"
$ws.Range("J11").Value = "tmp2e4a1"
$ws.Range("J12").Value = "tmp2e4a1"
$ws.Range("J13").Value = "tmp2e4a1"
$ws.Range("J14").Value = "tmp2e4a1"
$ws.Range("J15").Value = "tmp2e4a1"
$ws.Range("J16").Value = "tmp2e4a1"
$ws.Range("J17").Value = "tmp2e4a1"
$ws.Range("J18").Value = "tmp2e4a1"
$ws.Range("J19").Value = "tmp2e4a1"
$ws.Range("J20").Value = "tmp2e4a1"
$ws.Range("J21").Value = "tmp2e4a1"
$ws.Range("J22").Value = "tmp2e4a1"
$ws.Range("J23").Value = "tmp2e4a1"
$ws.Range("J24").Value = "tmp2e4a1"
$ws.Range("J25").Value = "tmp2e4a1"
$ws.Range("J26").Value = "tmp2e4a1"
$ws.Range("J27").Value = "tmp2e4a1"
$ws.Range("J28").Value = "tmp2e4a1"
$ws.Range("J29").Value = "tmp2e4a1"
$ws.Range("J30").Value = "tmp2e4a1"

# Column K: index (header + data)
$ws.Range("K1").Value = "index"
$ws.Range("K2").Value = 95
$ws.Range("K3").Value = 96
$ws.Range("K4").Value = 97
$ws.Range("K5").Value = 98
$ws.Range("K6").Value = 99
$ws.Range("K7").Value = 100
$ws.Range("K8").Value = 101
$ws.Range("K9").Value = 102
$ws.Range("K10").Value = 103
$ws.Range("K11").Value = 104
$ws.Range("K12").Value = 105
$ws.Range("K13").Value = 106
$ws.Range("K14").Value = 107
$ws.Range("K15").Value = 108
$ws.Range("K16").Value = 110
$ws.Range("K17").Value = 111
$ws.Range("K18").Value = 112
$ws.Range("K19").Value = 113
$ws.Range("K20").Value = 114
$ws.Range("K21").Value = 115
$ws.Range("K22").Value = 116
$ws.Range("K23").Value = 117
$ws.Range("K24").Value = 118
$ws.Range("K25").Value = 119
$ws.Range("K26").Value = 120
$ws.Range("K27").Value = 121
$ws.Range("K28").Value = 122
$ws.Range("K29").Value = 123
$ws.Range("K30").Value = 124

